$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.333.92"
$ws.Range("D3").Value = "'3.838.08"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'601.29"
$ws.Range("D6").Value = "'169.20"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("D7").Value = "'3.837.89"
$ws.Range("E7").Value = "  -1.49%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("D10").Value = "'0.166"
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("D11").Value = "'6.49"
$ws.Range("E11").Value = "  +1.88%  "
$ws.Range("D12").Value = "'0.458"
$ws.Range("E12").Value = "  -2.10%  "
$ws.Range("D13").Value = "'0.0000272"
$ws.Range("E13").Value = "  +5.93%  "
$ws.Range("D14").Value = "'37.12"
$ws.Range("E14").Value = "  -3.09%  "
$ws.Range("D15").Value = "'4.477.89"
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("D16").Value = "'3.827.41"
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("D17").Value = "'68.343.22"
$ws.Range("E17").Value = "  -1.94%  "
$ws.Range("D18").Value = "'18.55"
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("D19").Value = "'7.39"
$ws.Range("E19").Value = "  -3.48%  "
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("D21").Value = "'11.11"
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("D22").Value = "'470.30"
$ws.Range("E22").Value = "  -4.19%  "
$ws.Range("D23").Value = "'0.735"
$ws.Range("E23").Value = "  -1.93%  "
$ws.Range("D24").Value = "'0.0000159"
$ws.Range("E24").Value = "  -4.31%  "
$ws.Range("D25").Value = "'83.19"
$ws.Range("E25").Value = "  -2.58%  "
$ws.Range("E26").Value = "  -3.81%  "
$ws.Range("D27").Value = "'12.15"
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("D28").Value = "'10.06"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("D31").Value = "'3.984.23"
$ws.Range("E31").Value = "  -1.59%  "
$ws.Range("D32").Value = "'7.70"
$ws.Range("E32").Value = "  -1.67%  "
$ws.Range("D33").Value = "'31.63"
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("E34").Value = "  -4.68%  "
$ws.Range("D35").Value = "'9.39"
$ws.Range("E35").Value = "  -2.79%  "
$ws.Range("D36").Value = "'3.799.70"
$ws.Range("E36").Value = "  -1.70%  "
$ws.Range("E37").Value = "  -2.15%  "
$ws.Range("D38").Value = "'3.69"
$ws.Range("E38").Value = "  +10.22%  "
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").Value = "'1.02"
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.140"
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("D41").Value = "'5.95"
$ws.Range("E41").Value = "  -3.02%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").Value = "'0.315"
$ws.Range("E43").Value = "  -4.24%  "
$ws.Range("E44").Value = "  -5.73%  "
$ws.Range("D45").Value = "'8.75"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "'416.55"
$ws.Range("E47").Value = "  -4.79%  "
$ws.Range("D48").Value = "'47.17"
$ws.Range("E48").Value = "  -2.32%  "
$ws.Range("D49").Value = "'0.000287"
$ws.Range("E49").Value = "  +4.44%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0360"
$ws.Range("E50").Value = "  -2.44%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'141.88"
$ws.Range("E51").Value = "  -1.44%  "
